$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: heading text "ATRIBUIÇÕES" -> "RESPONSABILIDADE POR FUNÇÃO"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("ATRIBUIÇÕES", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "RESPONSABILIDADE POR FUNÇÃO", 2)

# ---------------------------------------------------------------------------
# Change 2: rewrite the paragraph that used to read
#   "Cabe a Chefe da NTDM garantir o cumprimento dos seguintes processos:"
# into
#   "As responsabilidades por função são delineadas sinteticamente no
#    Regimento Interno do CELOG e as atividades relacionados aos seus
#    cumprimentos estão contempladas nos seguintes processos:"
# while keeping the leading three tab runs intact and giving the first two
# new runs a szCs=24 (complex-script size) that the originals lacked.
# ---------------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("Cabe a Chefe da NTDM garantir o cumprimento dos seguintes processos:", `
                                $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Cabe a Chefe da NTDM...' paragraph text"
}

# Grab the whole owning paragraph, wipe its text (the pilcrow stays put),
# then rebuild it from scratch so we get exact control over every run.
$targetPara = $findRng.Paragraphs(1)
$paraRange = $targetPara.Range
$textOnly = $d.Range($paraRange.Start, $paraRange.End - 1)
$textOnly.Delete()

$insertAt = $d.Range($paraRange.Start, $paraRange.Start)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:spacing w:before="120"/><w:rPr><w:sz w:val="24"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="pt-PT"/></w:rPr><w:tab/></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="pt-PT"/></w:rPr><w:tab/></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="pt-PT"/></w:rPr><w:tab/></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t xml:space="preserve">As </w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>responsabilidades por função são delineadas sinteticamente no Regimento Interno do CELOG e as atividades relacionados aos seus cumprimentos estão contempladas n</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>os seguintes processos:</w:t></w:r>' + `
  '</w:p>'
$insertAt.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# Change 3: after the "... Gestão de calibração de equipamentos da NNAQ"
# bullet, add a blank underline-styled paragraph, a new level-1 heading
# "RESPONSABILIDADE NO INTER-RELACIONAMENTO ENTRE OS SETORES" and a new
# level-2 body paragraph describing it. These land just before the
# pre-existing blank paragraph that precedes "DISPOSIÇÕES FINAIS".
# ---------------------------------------------------------------------------
$anchorRng = $d.Content
$anchorFound = $anchorRng.Find.Execute("Gestão de calibração de equipamentos da NNAQ", `
                                        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $anchorFound) {
    throw "Could not locate the 'Gestão de calibração de equipamentos da NNAQ' paragraph"
}
$anchorPara = $anchorRng.Paragraphs(1)
# Insert one character before the paragraph's own pilcrow (not at the
# following paragraph's absolute start) - this engine merges a trailing
# <w:p> fragment into an already-existing "next" paragraph when the
# insertion point sits exactly at that following paragraph's start, which
# would clobber the pre-existing blank paragraph before "DISPOSIÇÕES
# FINAIS". Anchoring just inside the preceding paragraph avoids that.
$insertPos = $anchorPara.Range.End - 1
$insertRng = $d.Range($insertPos, $insertPos)

$blockXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/><w:spacing w:before="120"/><w:jc w:val="both"/>' + `
  '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr>' + `
  '</w:p>' + `
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/>' + `
  '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' + `
  '<w:spacing w:before="120"/><w:jc w:val="both"/>' + `
  '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr>' + `
  '<w:t>RESPONSABILIDADE NO INTER-RELACIONAMENTO ENTRE OS SETORES</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:keepNext/><w:keepLines/><w:widowControl w:val="0"/>' + `
  '<w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' + `
  '<w:spacing w:before="120"/><w:jc w:val="both"/>' + `
  '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="pt-PT"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr><w:t>As responsa</w:t></w:r>' + `
  '<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="pt-PT"/></w:rPr>' + `
  '<w:t>bilidades no inter-relacionamento entre setores são apresentadas detalhadamente nos PLOG relacionados no item 2.2 desta NPA.</w:t></w:r>' + `
  '</w:p>'
$insertRng.InsertXML($blockXml)

Write-Output "Edits applied"
